$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112, shifting existing rows 112-125 down to 113-126
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new weekly data point
$ws.Cells.Item(112, 1).Value = 5
$ws.Cells.Item(112, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(112, 3).Value = "Maule"
$ws.Cells.Item(112, 4).Value = 45142
$ws.Cells.Item(112, 5).Value = 7
$ws.Cells.Item(112, 6).Value = 100112013
$ws.Cells.Item(112, 7).Value = "Alcachofa"
$ws.Cells.Item(112, 8).Value = "Madrigal"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 300
$ws.Cells.Item(112, 11).Value = 14000
$ws.Cells.Item(112, 12).Value = 14000
$ws.Cells.Item(112, 13).Value = 14000
$ws.Cells.Item(112, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(112, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(112, 16).Value = 350
$ws.Cells.Item(112, 17).Value = 40
$ws.Cells.Item(112, 18).Value = "Hortaliza"
